$p = $ppt.ActivePresentation

# The deck currently carries two theme parts:
#   theme1.xml ("Integral")     -> used by the (only reachable) SlideMaster/Theme
#   theme2.xml ("Office Theme") -> used by the Notes Master
# The authored edit swaps their contents, so the slide master (and therefore
# the slides themselves) end up using the default "Office Theme" palette.
# Reproduce that swap on the one Theme object the COM surface exposes
# (SlideMaster.Theme), driving its ThemeColorScheme to the "Office Theme"
# colour values (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - in that
# fixed Item() order).

$officeThemeColors = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
